$wb = $excel.ActiveWorkbook

# --- Precondiciones sheet ---
$wsPre = $wb.Worksheets.Item("Precondiciones")

# Update existing precondition text (row 2) to reference "Córdoba" instead of <Ciudad1>
$wsPre.Range("B2").Value = 'Hay playas de estacionamiento cargadas en la Base de Datos para la ciudad "Córdoba"'

# Add a new precondition row (row 3) with matching formatting copied from row 2
$wsPre.Range("A3").Value = 2
$wsPre.Range("B3").Value = 'La ciudad ¨Códoba" existe en la base de datos.'
$wsPre.Range("A2:B2").Copy()
$wsPre.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPre.Rows.Item(2).RowHeight = 26.25

$wsPre.Range("B4").Select()

# --- Pasos sheet ---
$wsPasos = $wb.Worksheets.Item("Pasos")
$wsPasos.Range("B3").Value = 'Ingreso "Córdoba" en el campo nombre de ciudad'
$wsPasos.Range("C7").Select()

# --- Activate Precondiciones as the selected tab ---
$wsPre.Activate()
